$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.937.88'
$ws.Range("E2").Value = '  -1.96%  '
$ws.Range("D3").Value = '2.913.64'
$ws.Range("E3").Value = '  -3.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.03'
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.47'
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '2.912.82'
$ws.Range("E8").Value = '  -3.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("E9").Value = '  -2.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.71'
$ws.Range("E10").Value = '  +7.65%  '
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("E12").Value = '  -2.08%  '
$ws.Range("E13").Value = '  -2.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.36'
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = '3.397.03'
$ws.Range("E16").Value = '  -3.18%  '
$ws.Range("D17").Value = '60.922.46'
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.82'
$ws.Range("E18").Value = '  -2.48%  '
$ws.Range("D19").Value = '2.913.63'
$ws.Range("E19").Value = '  -3.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '426.03'
$ws.Range("E20").Value = '  -5.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.65'
$ws.Range("E21").Value = '  -3.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.670'
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.16'
$ws.Range("E23").Value = '  -2.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.84'
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.00'
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  -2.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.79'
$ws.Range("E27").Value = '  -2.60%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.23'
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("E31").Value = '  -2.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.17'
$ws.Range("E32").Value = '  +3.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.69'
$ws.Range("E33").Value = '  -2.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.107'
$ws.Range("E34").Value = '  -3.76%  '
$ws.Range("D35").Value = '0.0₃0846'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.67'
$ws.Range("E37").Value = '  -2.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.02'
$ws.Range("E38").Value = '  +3.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.82'
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E42").Value = '  -2.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.75'
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.287'
$ws.Range("E44").Value = '  +2.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '375.78'
$ws.Range("E45").Value = '  -5.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0345'
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("D47").Value = '2.653.29'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.90'
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.45'
$ws.Range("E49").Value = '  +7.42%  '
$ws.Range("E51").Value = '  -0.87%  '

Write-Host "Applied all cell updates"